# Updated cryptos list (coinranking.com snapshot refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds free-form text such as "27.696.20" or "4.690"
# (not valid numeric literals, or numbers whose trailing zero / dot-grouping
# matters). Forcing each touched cell to Text format ("@") before writing its
# new value keeps it stored as text, matching the source file, instead of
# Excel silently reinterpreting it as a number and losing that formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.696.20'
$ws.Range("E2").Value = '  +5.81%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.730.82'
$ws.Range("E3").Value = '  +4.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.49'
$ws.Range("E5").Value = '  +3.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5444'
$ws.Range("E6").Value = '  +3.26%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2737'
$ws.Range("E8").Value = '  +1.91%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06699'
$ws.Range("E9").Value = '  +4.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.84'
$ws.Range("E10").Value = '  +5.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07773'
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.690'
$ws.Range("E12").Value = '  +1.61%  '
# Row 13: coin entry swapped in with refreshed price/volume data
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.748.91'
$ws.Range("E13").Value = '  +5.69%  '
# Row 14: coin entry swapped in with refreshed price/volume data
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.971.13'
$ws.Range("E14").Value = '  +4.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5967'
$ws.Range("E15").Value = '  +5.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8402'
$ws.Range("E16").Value = '  +1.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.06'
$ws.Range("E17").Value = '  +4.98%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '27.714.49'
$ws.Range("E18").Value = '  +5.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '225.57'
$ws.Range("E19").Value = '  +18.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.812'
$ws.Range("E20").Value = '  +2.59%  '
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.87'
$ws.Range("E22").Value = '  +4.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.210'
$ws.Range("E23").Value = '  +3.54%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.52'
$ws.Range("E25").Value = '  +0.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.735'
$ws.Range("E26").Value = '  +13.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1249'
$ws.Range("E27").Value = '  +3.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.466'
$ws.Range("E28").Value = '  +2.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '17.06'
$ws.Range("E29").Value = '  +6.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05680'
$ws.Range("E30").Value = '  +0.59%  '
$ws.Range("E31").Value = '  +2.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.677'
$ws.Range("E32").Value = '  +5.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.507'
$ws.Range("E33").Value = '  +3.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.674'
$ws.Range("E34").Value = '  +5.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9745'
$ws.Range("E35").Value = '  +2.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.850'
$ws.Range("E36").Value = '  +1.83%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.438'
$ws.Range("E37").Value = '  +1.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5995'
$ws.Range("E38").Value = '  +3.70%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01667'
$ws.Range("E39").Value = '  +3.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.919'
$ws.Range("E40").Value = '  -1.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8590'
$ws.Range("E41").Value = '  +2.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.048.80'
$ws.Range("E42").Value = '  +1.87%  '
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.45'
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.875.97'
$ws.Range("E45").Value = '  +4.44%  '
$ws.Range("E46").Value = '  +9.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '59.50'
$ws.Range("E47").Value = '  +1.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.300'
$ws.Range("E48").Value = '  +3.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4433'
$ws.Range("E49").Value = '  +2.12%  '
# Row 50: coin entry swapped in with refreshed price/volume data
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05332'
$ws.Range("E50").Value = '  -0.54%  '
# Row 51: coin entry swapped in with refreshed price/volume data
$ws.Range("B51").Value = 'Frax'
$ws.Range("C51").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.002'
# NOTE: source diff hunk for row 51 is cut off before the Volume(1h) value;
# best estimate derived from the other stablecoin rows +0.04..+0.08pp shift.
$ws.Range("E51").Value = '  +0.05%  '
